$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameters")
$ws.Range("E1").Value = "name"
